# Fixed naive component forecaster bug - Presentation state 11.02.
# Updates the naive-forecast QoQ error values (columns B:K, rows 24-52)
# on Sheet1 to reflect the corrected forecaster output, including filling
# in previously-blank forecast-error cells for later vintages (rows 33-52).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(24, 11).Value = -8.183026219731815
$ws.Cells.Item(25, 10).Value = -8.028721684813732
$ws.Cells.Item(25, 11).Value = 0.1109553426839227
$ws.Cells.Item(26, 9).Value = -8.680305474188014
$ws.Cells.Item(26, 10).Value = -0.54062844669036
$ws.Cells.Item(26, 11).Value = 1.60694255667963
$ws.Cells.Item(27, 8).Value = -8.462984671940774
$ws.Cells.Item(27, 9).Value = -0.3233076444431192
$ws.Cells.Item(27, 10).Value = 1.824263358926871
$ws.Cells.Item(27, 11).Value = -1.615771047272361
$ws.Cells.Item(28, 7).Value = -8.080305474188014
$ws.Cells.Item(28, 8).Value = 0.05937155330964
$ws.Cells.Item(28, 9).Value = 2.20694255667963
$ws.Cells.Item(28, 10).Value = -1.233091849519602
$ws.Cells.Item(28, 11).Value = -1.294275412034003
$ws.Cells.Item(29, 6).Value = -8.554988394376949
$ws.Cells.Item(29, 7).Value = -0.4153113668792941
$ws.Cells.Item(29, 8).Value = 1.732259636490696
$ws.Cells.Item(29, 9).Value = -1.707774769708536
$ws.Cells.Item(29, 10).Value = -1.768958332222937
$ws.Cells.Item(29, 11).Value = 0.2725059170203678
$ws.Cells.Item(30, 5).Value = -8.396144372968052
$ws.Cells.Item(30, 6).Value = -0.2564673454703979
$ws.Cells.Item(30, 7).Value = 1.891103657899592
$ws.Cells.Item(30, 8).Value = -1.54893074829964
$ws.Cells.Item(30, 9).Value = -1.610114310814041
$ws.Cells.Item(30, 10).Value = 0.431349938429264
$ws.Cells.Item(30, 11).Value = -0.1512409947814746
$ws.Cells.Item(31, 4).Value = -8.452313642352863
$ws.Cells.Item(31, 5).Value = -0.3126366148552081
$ws.Cells.Item(31, 6).Value = 1.834934388514782
$ws.Cells.Item(31, 7).Value = -1.60510001768445
$ws.Cells.Item(31, 8).Value = -1.666283580198851
$ws.Cells.Item(31, 9).Value = 0.3751806690444539
$ws.Cells.Item(31, 10).Value = -0.2074102641662847
$ws.Cells.Item(31, 11).Value = -0.1117330153097526
$ws.Cells.Item(32, 3).Value = -10.68030547418801
$ws.Cells.Item(32, 4).Value = -2.54062844669036
$ws.Cells.Item(32, 5).Value = -0.3930574433203702
$ws.Cells.Item(32, 6).Value = -3.833091849519602
$ws.Cells.Item(32, 7).Value = -3.894275412034003
$ws.Cells.Item(32, 8).Value = -1.852811162790698
$ws.Cells.Item(32, 9).Value = -2.435402096001437
$ws.Cells.Item(32, 10).Value = -2.339724847144905
$ws.Cells.Item(32, 11).Value = -2.60014685520589
$ws.Cells.Item(33, 2).Value = -18.17126098581633
$ws.Cells.Item(33, 3).Value = -10.03158395831868
$ws.Cells.Item(33, 4).Value = -7.884012954948687
$ws.Cells.Item(33, 5).Value = -11.32404736114792
$ws.Cells.Item(33, 6).Value = -11.38523092366232
$ws.Cells.Item(33, 7).Value = -9.343766674419015
$ws.Cells.Item(33, 8).Value = -9.926357607629754
$ws.Cells.Item(33, 9).Value = -9.830680358773222
$ws.Cells.Item(33, 10).Value = -10.09110236683421
$ws.Cells.Item(33, 11).Value = -9.255644381000664
$ws.Cells.Item(34, 2).Value = 8.139677027497655
$ws.Cells.Item(34, 3).Value = 10.28724803086764
$ws.Cells.Item(34, 4).Value = 6.847213624668413
$ws.Cells.Item(34, 5).Value = 6.786030062154012
$ws.Cells.Item(34, 6).Value = 8.827494311397317
$ws.Cells.Item(34, 7).Value = 8.244903378186578
$ws.Cells.Item(34, 8).Value = 8.34058062704311
$ws.Cells.Item(34, 9).Value = 8.080158618982125
$ws.Cells.Item(34, 10).Value = 8.915616604815668
$ws.Cells.Item(34, 11).Value = 8.81559383571846
$ws.Cells.Item(35, 2).Value = 2.14757100336999
$ws.Cells.Item(35, 3).Value = -1.292463402829242
$ws.Cells.Item(35, 4).Value = -1.353646965343643
$ws.Cells.Item(35, 5).Value = 0.6878172838996619
$ws.Cells.Item(35, 6).Value = 0.1052263506889233
$ws.Cells.Item(35, 7).Value = 0.2009035995454554
$ws.Cells.Item(35, 8).Value = -0.05951840851552959
$ws.Cells.Item(35, 9).Value = 0.7759395773180131
$ws.Cells.Item(35, 10).Value = 0.6759168082208049
$ws.Cells.Item(35, 11).Value = 0.3220587481868762
$ws.Cells.Item(36, 2).Value = -3.440034406199232
$ws.Cells.Item(36, 3).Value = -3.501217968713633
$ws.Cells.Item(36, 4).Value = -1.459753719470328
$ws.Cells.Item(36, 5).Value = -2.042344652681066
$ws.Cells.Item(36, 6).Value = -1.946667403824535
$ws.Cells.Item(36, 7).Value = -2.20708941188552
$ws.Cells.Item(36, 8).Value = -1.371631426051977
$ws.Cells.Item(36, 9).Value = -1.471654195149185
$ws.Cells.Item(36, 10).Value = -1.825512255183114
$ws.Cells.Item(36, 11).Value = -1.677319185747749
$ws.Cells.Item(37, 2).Value = -0.06118356251440082
$ws.Cells.Item(37, 3).Value = 1.980280686728904
$ws.Cells.Item(37, 4).Value = 1.397689753518165
$ws.Cells.Item(37, 5).Value = 1.493367002374697
$ws.Cells.Item(37, 6).Value = 1.232944994313713
$ws.Cells.Item(37, 7).Value = 2.068402980147255
$ws.Cells.Item(37, 8).Value = 1.968380211050047
$ws.Cells.Item(37, 9).Value = 1.614522151016118
$ws.Cells.Item(37, 10).Value = 1.762715220451483
$ws.Cells.Item(37, 11).Value = 1.920148881959748
$ws.Cells.Item(38, 2).Value = 2.041464249243305
$ws.Cells.Item(38, 3).Value = 1.458873316032566
$ws.Cells.Item(38, 4).Value = 1.554550564889098
$ws.Cells.Item(38, 5).Value = 1.294128556828113
$ws.Cells.Item(38, 6).Value = 2.129586542661656
$ws.Cells.Item(38, 7).Value = 2.029563773564448
$ws.Cells.Item(38, 8).Value = 1.675705713530519
$ws.Cells.Item(38, 9).Value = 1.823898782965884
$ws.Cells.Item(38, 10).Value = 1.981332444474148
$ws.Cells.Item(38, 11).Value = 1.480676358197826
$ws.Cells.Item(39, 2).Value = -0.5825909332107386
$ws.Cells.Item(39, 3).Value = -0.4869136843542065
$ws.Cells.Item(39, 4).Value = -0.7473356924151915
$ws.Cells.Item(39, 5).Value = 0.08812229341835121
$ws.Cells.Item(39, 6).Value = -0.01190047567885699
$ws.Cells.Item(39, 7).Value = -0.3657585357127857
$ws.Cells.Item(39, 8).Value = -0.2175654662774206
$ws.Cells.Item(39, 9).Value = -0.06013180476915631
$ws.Cells.Item(39, 10).Value = -0.5607878910454785
$ws.Cells.Item(39, 11).Value = -0.2804316037748719
$ws.Cells.Item(40, 2).Value = 0.09567724885653209
$ws.Cells.Item(40, 3).Value = -0.1647447592044529
$ws.Cells.Item(40, 4).Value = 0.6707132266290898
$ws.Cells.Item(40, 5).Value = 0.5706904575318816
$ws.Cells.Item(40, 6).Value = 0.2168323974979529
$ws.Cells.Item(40, 7).Value = 0.365025466933318
$ws.Cells.Item(40, 8).Value = 0.5224591284415823
$ws.Cells.Item(40, 9).Value = 0.02180304216525999
$ws.Cells.Item(40, 10).Value = 0.3021593294358667
$ws.Cells.Item(40, 11).Value = 0.1301932178764815
$ws.Cells.Item(41, 2).Value = -0.260422008060985
$ws.Cells.Item(41, 3).Value = 0.5750359777725577
$ws.Cells.Item(41, 4).Value = 0.4750132086753495
$ws.Cells.Item(41, 5).Value = 0.1211551486414208
$ws.Cells.Item(41, 6).Value = 0.2693482180767859
$ws.Cells.Item(41, 7).Value = 0.4267818795850502
$ws.Cells.Item(41, 8).Value = -0.07387420669127209
$ws.Cells.Item(41, 9).Value = 0.2064820805793346
$ws.Cells.Item(41, 10).Value = 0.0345159690199494
$ws.Cells.Item(41, 11).Value = 0.3402630873498395
$ws.Cells.Item(42, 2).Value = 0.8354579858335427
$ws.Cells.Item(42, 3).Value = 0.7354352167363345
$ws.Cells.Item(42, 4).Value = 0.3815771567024058
$ws.Cells.Item(42, 5).Value = 0.5297702261377709
$ws.Cells.Item(42, 6).Value = 0.6872038876460351
$ws.Cells.Item(42, 7).Value = 0.1865478013697129
$ws.Cells.Item(42, 8).Value = 0.4669040886403196
$ws.Cells.Item(42, 9).Value = 0.2949379770809344
$ws.Cells.Item(42, 10).Value = 0.6006850954108245
$ws.Cells.Item(42, 11).Value = -0.01150507022163078
$ws.Cells.Item(43, 2).Value = -0.1000227690972082
$ws.Cells.Item(43, 3).Value = -0.4538808291311369
$ws.Cells.Item(43, 4).Value = -0.3056877596957718
$ws.Cells.Item(43, 5).Value = -0.1482540981875075
$ws.Cells.Item(43, 6).Value = -0.6489101844638298
$ws.Cells.Item(43, 7).Value = -0.3685538971932231
$ws.Cells.Item(43, 8).Value = -0.5405200087526083
$ws.Cells.Item(43, 9).Value = -0.2347728904227182
$ws.Cells.Item(43, 10).Value = -0.8469630560551735
$ws.Cells.Item(43, 11).Value = -0.1586268237156929
$ws.Cells.Item(44, 2).Value = -0.3538580600339287
$ws.Cells.Item(44, 3).Value = -0.2056649905985636
$ws.Cells.Item(44, 4).Value = -0.04823132909029931
$ws.Cells.Item(44, 5).Value = -0.5488874153666216
$ws.Cells.Item(44, 6).Value = -0.2685311280960149
$ws.Cells.Item(44, 7).Value = -0.4404972396554001
$ws.Cells.Item(44, 8).Value = -0.13475012132551
$ws.Cells.Item(44, 9).Value = -0.7469402869579653
$ws.Cells.Item(44, 10).Value = -0.05860405461848467
$ws.Cells.Item(45, 2).Value = 0.1481930694353651
$ws.Cells.Item(45, 3).Value = 0.3056267309436294
$ws.Cells.Item(45, 4).Value = -0.1950293553326929
$ws.Cells.Item(45, 5).Value = 0.08532693193791374
$ws.Cells.Item(45, 6).Value = -0.08663917962147143
$ws.Cells.Item(45, 7).Value = 0.2191079387084187
$ws.Cells.Item(45, 8).Value = -0.3930822269240366
$ws.Cells.Item(45, 9).Value = 0.295254005415444
$ws.Cells.Item(46, 2).Value = 0.1574336615082643
$ws.Cells.Item(46, 3).Value = -0.343222424768058
$ws.Cells.Item(46, 4).Value = -0.06286613749745135
$ws.Cells.Item(46, 5).Value = -0.2348322490568365
$ws.Cells.Item(46, 6).Value = 0.0709148692730536
$ws.Cells.Item(46, 7).Value = -0.5412752963594016
$ws.Cells.Item(46, 8).Value = 0.1470609359800789
$ws.Cells.Item(47, 2).Value = -0.5006560862763223
$ws.Cells.Item(47, 3).Value = -0.2202997990057156
$ws.Cells.Item(47, 4).Value = -0.3922659105651008
$ws.Cells.Item(47, 5).Value = -0.08651879223521067
$ws.Cells.Item(47, 6).Value = -0.698708957867666
$ws.Cells.Item(47, 7).Value = -0.01037272552818536
$ws.Cells.Item(48, 2).Value = 0.2803562872706067
$ws.Cells.Item(48, 3).Value = 0.1083901757112215
$ws.Cells.Item(48, 4).Value = 0.4141372940411116
$ws.Cells.Item(48, 5).Value = -0.1980528715913437
$ws.Cells.Item(48, 6).Value = 0.4902833607481369
$ws.Cells.Item(49, 2).Value = -0.1719661115593852
$ws.Cells.Item(49, 3).Value = 0.1337810067705049
$ws.Cells.Item(49, 4).Value = -0.4784091588619503
$ws.Cells.Item(49, 5).Value = 0.2099270734775303
$ws.Cells.Item(50, 2).Value = 0.3057471183298901
$ws.Cells.Item(50, 3).Value = -0.3064430473025652
$ws.Cells.Item(50, 4).Value = 0.3818931850369154
$ws.Cells.Item(51, 2).Value = -0.6121901656324553
$ws.Cells.Item(51, 3).Value = 0.07614606670702531
$ws.Cells.Item(52, 2).Value = 0.6883362323394806
